$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold/centered) from I1 into the two new header cells, then set their text
$ws.Range("I1").Copy($ws.Range("J1:K1"))
$ws.Range("J1").Value = 'Onkelos'
$ws.Range("K1").Value = 'Jonathan'

# Copy body style (wrap text) from I2/I3 into the new data cells, then set their text
$ws.Range("I2").Copy($ws.Range("J2:K2"))
$ws.Range("J2").Value = 'And there has not ever arisen a prophet within Yisroel like Moshe, whom Adonoy knew [<b>appeared to</b>] face-to-face.'
$ws.Range("K2").Value = 'But no prophet hath again risen in Israel like unto Mosheh, because the Word of the Lord had known him to speak with him word for word,'

$ws.Range("I3").Copy($ws.Range("J3:K3"))
$ws.Range("J3").Value = '“Go, gather the elders of Yisrael, and say to them, ‘Adonoy, the God of your fathers appeared [<b>became revealed</b>] to me—the God of Avraham, Yitzchok and Yaakov—saying, “I have indeed been mindful of you, regarding that which is being done to you in Egypt.'
$ws.Range("K3").Value = 'Go, and assemble the elders of Israel, and say to them, The Lord God of your fathers hath appeared unto me, the God of Abraham, Izhak, and Jakob, saying, Remembering, I have remembered you, and the injury that is done you in Mizraim;'

# Set the new column widths (OOXML stored width = ColumnWidth + 5/6)
$ws.Columns.Item(10).ColumnWidth = 313.56666666666666
$ws.Columns.Item(11).ColumnWidth = 281.16666666666669

